$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use single-quoted (literal, non-interpolating) strings so the '$' currency
# signs in these values are written verbatim.
$ws.Range("A2").Value = '268 444 6 86 8 26 86 8 28 6 8 2K 2'
$ws.Range("D2").Value = '$ 386 , 523 , 434 , 322.15 AND 26 %'
$ws.Range("G2").Value = '6E 8 2 4 8  2 6 8 6 2 6 8'
$ws.Range("H2").Value = '$ 286 , 027 , 341 , 398.39 AND $ 1 , 305 , 350 , 268.83'
$ws.Range("I2").Value = '$ 0.00 AND $ 42 , 595 , 191 , 764.52'
$ws.Range("J2").Value = '$ 76 , 273 , 957.85 AND NA'
